# Weekly fruit/vegetable price update: insert a new daily record at row 127
# (shifting the existing rows 127-158 down to 128-159) for the
# "Vega Modelo de Temuco - Pepino dulce" subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 127; this pushes the old rows 127:158 down
# to 128:159 (and extends the sheet dimension from R158 to R159).
$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new weekly record.
$ws.Cells.Item(127, 1).Value  = 10
$ws.Cells.Item(127, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value  = "La Araucanía"
$ws.Cells.Item(127, 4).Value  = 44476
$ws.Cells.Item(127, 5).Value  = 9
$ws.Cells.Item(127, 6).Value  = 100112043
$ws.Cells.Item(127, 7).Value  = "Pepino dulce"
$ws.Cells.Item(127, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(127, 9).Value  = "Segunda"
$ws.Cells.Item(127, 10).Value = 80
$ws.Cells.Item(127, 11).Value = 20000
$ws.Cells.Item(127, 12).Value = 20000
$ws.Cells.Item(127, 13).Value = 20000
$ws.Cells.Item(127, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(127, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(127, 16).Value = 1111
$ws.Cells.Item(127, 17).Value = 18
$ws.Cells.Item(127, 18).Value = "Hortaliza"
